$wb = $excel.ActiveWorkbook

# Add the new shared string value "some data" into a single-cell A1 of
# the sheets that were previously skipped because they only had a single
# part (bug: "Fixed skipping of single part sets").
$sheetNames = @("partInParts", "setPart1", "setPart2", "setPart3")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A1").Value = "some data"
    $ws.Range("A1").Style = "Normal"
}

# Fix the active sheet / tab selection: previously sheet index 8
# (setPart3) was marked active/selected; now sheet index 2 (partInParts)
# should be the active tab.
$wb.Worksheets.Item("partInParts").Activate()
